$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "THIS AGREEMENT is made at 02:18 on this 15 day of JUNE, 2023" "THIS AGREEMENT is made at 13:40 on this 15 day of JUNE, 2020"

Replace-Text "has a paid up capital of Rs. #49 as on 15/JUNE/2023." "has a paid up capital of Rs. #49 as on 15/JUNE/2020."

Replace-Text "inter corporate deposit of Rs. 40000/- Rupees only for a period of 5000 days beginning from the date of disbursal of loan i.e. 15/JUNE/2023." "inter corporate deposit of Rs. 10000/- Rupees only for a period of 10 days beginning from the date of disbursal of loan i.e. 15/JUNE/2020."

Replace-Text "inter-corporate deposit of Rs.40000 Rupees only to the Borrower" "inter-corporate deposit of Rs.10000 Rupees only to the Borrower"

Replace-Text "inter corporate deposit of Rs. 40000 Rupees only for a period of 5000 days beginning from the date of disbursal." "inter corporate deposit of Rs. 10000 Rupees only for a period of 10 days beginning from the date of disbursal."

Replace-Text "shall carry an interest @ 2% per annum payable" "shall carry an interest @ 1% per annum payable"

Replace-Text "pay a penal interest @ 2% per annum over" "pay a penal interest @ 1% per annum over"

Replace-Text "listed at 4000 Stock Exchange and the current market price of shares is agreed to be Rs. 421/- per share." "listed at 100 Stock Exchange and the current market price of shares is agreed to be Rs. 100/- per share."

Replace-Text "so as to ensure 2% margin between" "so as to ensure 1% margin between"

Replace-Text "in respect of the said 4999 Equity Shares of Metro" "in respect of the said 500 Equity Shares of Metro"
